# Fill in the missing forecast-error values for the GVA ifoCAST error series.
# This mirrors the preprocessing step that appends the newest "diff" columns
# for rows 136-145 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D136").Value = 0.7254492243564907
$ws.Range("D137").Value = 0.7215746373564907
$ws.Range("D138").Value = 0.5311946523564907
$ws.Range("D139").Value = 0.5539812373564907
$ws.Range("C140").Value = 0.1753415943564907
$ws.Range("C141").Value = 0.2651053283564908
$ws.Range("C142").Value = 0.08763596535649075
$ws.Range("C143").Value = 0.1003532183564907
$ws.Range("C144").Value = -0.02418658464350926
$ws.Range("C145").Value = 0.2001520573564908
